# Changing art start date to EpisodeOfCare resource
# - Remove the "TPT Start Date" (tpt-start-date) and "TPT Stop Date"
#   (tpt-stop-date) concepts from the Concepts sheet.
# - Update the Metadata sheet's Date and Count fields accordingly.

$wb = $excel.ActiveWorkbook

$wsConcepts = $wb.Worksheets.Item("Concepts")

# Row 4 holds "tpt-start-date" / "TPT Start Date".
# Deleting it shifts row 5 ("tpt-stop-date" / "TPT Stop Date") up into row 4,
# so deleting row 4 a second time removes that one too.
$wsConcepts.Rows.Item(4).Delete()
$wsConcepts.Rows.Item(4).Delete()

$wsMeta = $wb.Worksheets.Item("Metadata")

# Refresh the generated timestamp and the concept count (13 -> 11).
$wsMeta.Range("B8").Value = "2025-09-10T10:05:20+00:00"
# Leading apostrophe keeps "11" stored as text (matching the source data,
# which keeps Count as a string), instead of Excel auto-converting it to a number.
$wsMeta.Range("B22").Value = "'11"
